# Automatic update of files.
# Increment the "Förändrad" (changed) date in column C for rows 2-5 by one day
# (from 2023-09-15 / serial 45184 to 2023-09-16 / serial 45185).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    $cell.Value = $current.AddDays(1)
}
